$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title and source citation (January 2017 -> February 2017) ---
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("B49").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Rewrite the Annual Production table (rows 27-46) with the new country order/values ---
$ws.Cells.Item(27, 2).Value = "United States"
$ws.Cells.Item(27, 3).Value = 15.123586068
$ws.Cells.Item(27, 4).Value = 14.837390362000001
$ws.Cells.Item(27, 5).Value = 15.167720308
$ws.Cells.Item(27, 6).Value = 16.164785083000002

$ws.Cells.Item(28, 2).Value = "Canada"
$ws.Cells.Item(28, 3).Value = 4.5059315615999997
$ws.Cells.Item(28, 4).Value = 4.5256855260000002
$ws.Cells.Item(28, 5).Value = 4.7106092711000001
$ws.Cells.Item(28, 6).Value = 4.8683554024999998

$ws.Cells.Item(29, 2).Value = "Russia"
$ws.Cells.Item(29, 3).Value = 11.029721986
$ws.Cells.Item(29, 4).Value = 11.240399442999999
$ws.Cells.Item(29, 5).Value = 11.296731380000001
$ws.Cells.Item(29, 6).Value = 11.381449769

$ws.Cells.Item(30, 2).Value = "Brazil"
$ws.Cells.Item(30, 3).Value = 3.1831543781999998
$ws.Cells.Item(30, 4).Value = 3.2348831847000001
$ws.Cells.Item(30, 5).Value = 3.2955568134000002
$ws.Cells.Item(30, 6).Value = 3.3553912459999999

$ws.Cells.Item(31, 2).Value = "Kazakhstan"
$ws.Cells.Item(31, 3).Value = 1.7515083425
$ws.Cells.Item(31, 4).Value = 1.731690647
$ws.Cells.Item(31, 5).Value = 1.8647197606000001
$ws.Cells.Item(31, 6).Value = 1.8775569644000001

$ws.Cells.Item(32, 2).Value = "Australia"
$ws.Cells.Item(32, 3).Value = 0.41624702740000002
$ws.Cells.Item(32, 4).Value = 0.39293790599
$ws.Cells.Item(32, 5).Value = 0.40674383335999997
$ws.Cells.Item(32, 6).Value = 0.46944442683999998

$ws.Cells.Item(33, 2).Value = "Oman"
$ws.Cells.Item(33, 3).Value = 0.99025742466
$ws.Cells.Item(33, 4).Value = 1.0181907205
$ws.Cells.Item(33, 5).Value = 1.0046012658000001
$ws.Cells.Item(33, 6).Value = 1.0332822758

$ws.Cells.Item(34, 2).Value = "Malaysia"
$ws.Cells.Item(34, 3).Value = 0.73517808219000003
$ws.Cells.Item(34, 4).Value = 0.74239094207
$ws.Cells.Item(34, 5).Value = 0.74508896690000004
$ws.Cells.Item(34, 6).Value = 0.74841519856000005

$ws.Cells.Item(35, 2).Value = "Syria"
$ws.Cells.Item(35, 3).Value = 0.034879
$ws.Cells.Item(35, 4).Value = 0.034493823010999999
$ws.Cells.Item(35, 5).Value = 0.031856231695000001
$ws.Cells.Item(35, 6).Value = 0.028652780209999999

$ws.Cells.Item(36, 2).Value = "India"
$ws.Cells.Item(36, 3).Value = 1.0097811739
$ws.Cells.Item(36, 4).Value = 0.99406415484999999
$ws.Cells.Item(36, 5).Value = 0.99876149835000005
$ws.Cells.Item(36, 6).Value = 0.99791301851000003

$ws.Cells.Item(37, 2).Value = "Norway"
$ws.Cells.Item(37, 3).Value = 1.9577091781
$ws.Cells.Item(37, 4).Value = 2.0064891313
$ws.Cells.Item(37, 5).Value = 1.9923885728999999
$ws.Cells.Item(37, 6).Value = 1.9211752388000001

$ws.Cells.Item(38, 2).Value = "Other North Sea"
$ws.Cells.Item(38, 3).Value = 0.17882209589
$ws.Cells.Item(38, 4).Value = 0.15476272671999999
$ws.Cells.Item(38, 5).Value = 0.14575620112000001
$ws.Cells.Item(38, 6).Value = 0.13431993077000001

$ws.Cells.Item(39, 2).Value = "Egypt"
$ws.Cells.Item(39, 3).Value = 0.70536719177999996
$ws.Cells.Item(39, 4).Value = 0.69097386152999996
$ws.Cells.Item(39, 5).Value = 0.67871741633000005
$ws.Cells.Item(39, 6).Value = 0.66011603700999999

$ws.Cells.Item(40, 2).Value = "Vietnam"
$ws.Cells.Item(40, 3).Value = 0.35269502740000003
$ws.Cells.Item(40, 4).Value = 0.32241535187999998
$ws.Cells.Item(40, 5).Value = 0.31649729056999998
$ws.Cells.Item(40, 6).Value = 0.30280423529

$ws.Cells.Item(41, 2).Value = "Sudan/S. Sudan"
$ws.Cells.Item(41, 3).Value = 0.26100813698999997
$ws.Cells.Item(41, 4).Value = 0.25739596085999999
$ws.Cells.Item(41, 5).Value = 0.25247352696000003
$ws.Cells.Item(41, 6).Value = 0.20252747583

$ws.Cells.Item(42, 2).Value = "United Kingdom"
$ws.Cells.Item(42, 3).Value = 0.92537488492999997
$ws.Cells.Item(42, 4).Value = 0.98198033359000003
$ws.Cells.Item(42, 5).Value = 0.91216140754999997
$ws.Cells.Item(42, 6).Value = 0.86584426287000005

$ws.Cells.Item(43, 2).Value = "Azerbaijan"
$ws.Cells.Item(43, 3).Value = 0.85822120000000002
$ws.Cells.Item(43, 4).Value = 0.84808225473999999
$ws.Cells.Item(43, 5).Value = 0.78740074248000003
$ws.Cells.Item(43, 6).Value = 0.76913306114000002

$ws.Cells.Item(44, 2).Value = "Colombia"
$ws.Cells.Item(44, 3).Value = 1.0292200718
$ws.Cells.Item(44, 4).Value = 0.92420220444000001
$ws.Cells.Item(44, 5).Value = 0.90841527583000004
$ws.Cells.Item(44, 6).Value = 0.90233059174999997

$ws.Cells.Item(45, 2).Value = "Mexico"
$ws.Cells.Item(45, 3).Value = 2.6246649178000001
$ws.Cells.Item(45, 4).Value = 2.4940562623
$ws.Cells.Item(45, 5).Value = 2.2793745109999999
$ws.Cells.Item(45, 6).Value = 2.2761747310999998

$ws.Cells.Item(46, 2).Value = "China"
$ws.Cells.Item(46, 3).Value = 5.1676419178000002
$ws.Cells.Item(46, 4).Value = 4.8555322403999996
$ws.Cells.Item(46, 5).Value = 4.7078556608
$ws.Cells.Item(46, 6).Value = 4.6496078618999999

# --- Update the Total Non-OPEC row (row 48) ---
$ws.Range("C48").Value = 58.795748158999999
$ws.Range("D48").Value = 58.197843460999998
$ws.Range("E48").Value = 58.479249150000001
$ws.Range("F48").Value = 59.547778907999998

# --- Add the new "K" total-growth column (=SUM of the three yearly deltas) for rows 27-46 ---
$ws.Range("K27:K46").Formula = "=SUM(H27:J27)"
